# chore: update Sheets via scheduled runner
# Refreshes cached marketboard-derived price/profit figures (columns H-N)
# across several leve rows on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(13, 8).Value = 22161  # H13: 27961 -> 22161
$ws.Cells.Item(13, 9).Value = 15251.25  # I13: 20001.666 -> 15251.25
$ws.Cells.Item(13, 10).Value = 49800  # J13: 39900 -> 49800
$ws.Cells.Item(13, 11).Value = 15251.25  # K13: 20001.666 -> 15251.25
$ws.Cells.Item(13, 12).Value = 49800  # L13: 39900 -> 49800
$ws.Cells.Item(13, 13).Value = -15082.25  # M13: -19832.666 -> -15082.25
$ws.Cells.Item(13, 14).Value = -50138  # N13: -40238 -> -50138

$ws.Cells.Item(15, 8).Value = 1432.836  # H15: 1259.8833 -> 1432.836
$ws.Cells.Item(15, 9).Value = 1432.836  # I15: 1259.8833 -> 1432.836
$ws.Cells.Item(15, 11).Value = 4298.508  # K15: 3779.6499 -> 4298.508
$ws.Cells.Item(15, 13).Value = -4129.508  # M15: -3610.6499 -> -4129.508

$ws.Cells.Item(51, 8).Value = 0  # H51: 5000 -> 0
$ws.Cells.Item(51, 9).Value = 0  # I51: 5000 -> 0
$ws.Cells.Item(51, 11).Value = 0  # K51: 5000 -> 0
$ws.Cells.Item(51, 13).ClearContents()  # M51: -4516 -> (removed)

$ws.Cells.Item(107, 8).Value = 534  # H107: 557.9524 -> 534
$ws.Cells.Item(107, 9).Value = 299.54544  # I107: 321.5 -> 299.54544
$ws.Cells.Item(107, 10).Value = 748.9167  # J107: 772.9091 -> 748.9167
$ws.Cells.Item(107, 11).Value = 299.54544  # K107: 321.5 -> 299.54544
$ws.Cells.Item(107, 12).Value = 748.9167  # L107: 772.9091 -> 748.9167
$ws.Cells.Item(107, 13).Value = 1620.45456  # M107: 1598.5 -> 1620.45456
$ws.Cells.Item(107, 14).Value = -4588.9167  # N107: -4612.9091 -> -4588.9167

$ws.Cells.Item(125, 8).Value = 143632.14  # H125: 67618 -> 143632.14
$ws.Cells.Item(125, 9).Value = 500399.5  # I125: 500515.5 -> 500399.5
$ws.Cells.Item(125, 10).Value = 925.2  # J125: 1018.38464 -> 925.2
$ws.Cells.Item(125, 11).Value = 4503595.5  # K125: 4504639.5 -> 4503595.5
$ws.Cells.Item(125, 12).Value = 8326.800000000001  # L125: 9165.46176 -> 8326.800000000001
$ws.Cells.Item(125, 13).Value = -4501135.5  # M125: -4502179.5 -> -4501135.5
$ws.Cells.Item(125, 14).Value = -13246.8  # N125: -14085.46176 -> -13246.8

$ws.Cells.Item(133, 8).Value = 42000  # H133: 0 -> 42000
$ws.Cells.Item(133, 10).Value = 42000  # J133: 0 -> 42000
$ws.Cells.Item(133, 12).Value = 42000  # L133: 0 -> 42000
$ws.Cells.Item(133, 14).Value = -52120  # N133: None -> -52120

$ws.Cells.Item(135, 8).Value = 1248.9412  # H135: 917.2353000000001 -> 1248.9412
$ws.Cells.Item(135, 9).Value = 1326  # I135: 854.72 -> 1326
$ws.Cells.Item(135, 10).Value = 1138.8572  # J135: 1090.8889 -> 1138.8572
$ws.Cells.Item(135, 11).Value = 11934  # K135: 7692.48 -> 11934
$ws.Cells.Item(135, 12).Value = 10249.7148  # L135: 9818.000099999999 -> 10249.7148
$ws.Cells.Item(135, 13).Value = -9399  # M135: -5157.48 -> -9399
$ws.Cells.Item(135, 14).Value = -15319.7148  # N135: -14888.0001 -> -15319.7148

$ws.Cells.Item(137, 8).Value = 1296.125  # H137: 1336.1333 -> 1296.125
$ws.Cells.Item(137, 9).Value = 1224.7826  # I137: 1215.875 -> 1224.7826
$ws.Cells.Item(137, 10).Value = 1478.4445  # J137: 1817.1666 -> 1478.4445
$ws.Cells.Item(137, 11).Value = 3674.3478  # K137: 3647.625 -> 3674.3478
$ws.Cells.Item(137, 12).Value = 4435.333500000001  # L137: 5451.4998 -> 4435.333500000001
$ws.Cells.Item(137, 13).Value = -1124.3478  # M137: -1097.625 -> -1124.3478
$ws.Cells.Item(137, 14).Value = -9535.333500000001  # N137: -10551.4998 -> -9535.333500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(76, 8).Value = 21969.715  # H76: 0 -> 21969.715
$ws.Cells.Item(76, 10).Value = 21969.715  # J76: 0 -> 21969.715
$ws.Cells.Item(76, 12).Value = 21969.715  # L76: 0 -> 21969.715
$ws.Cells.Item(76, 14).Value = -22645.715  # N76: None -> -22645.715

$ws.Cells.Item(79, 8).Value = 21969.715  # H79: 0 -> 21969.715
$ws.Cells.Item(79, 10).Value = 21969.715  # J79: 0 -> 21969.715
$ws.Cells.Item(79, 12).Value = 21969.715  # L79: 0 -> 21969.715
$ws.Cells.Item(79, 14).Value = -24309.715  # N79: None -> -24309.715

$ws.Cells.Item(92, 8).Value = 89622.86  # H92: 80045 -> 89622.86
$ws.Cells.Item(92, 10).Value = 89622.86  # J92: 80045 -> 89622.86
$ws.Cells.Item(92, 12).Value = 89622.86  # L92: 80045 -> 89622.86
$ws.Cells.Item(92, 14).Value = -94614.86  # N92: -85037 -> -94614.86

$ws.Cells.Item(122, 8).Value = 4623.7036  # H122: 5644.2104 -> 4623.7036
$ws.Cells.Item(122, 9).Value = 4953.933  # I122: 7462.375 -> 4953.933
$ws.Cells.Item(122, 10).Value = 4210.9165  # J122: 4321.909 -> 4210.9165
$ws.Cells.Item(122, 11).Value = 14861.799  # K122: 22387.125 -> 14861.799
$ws.Cells.Item(122, 12).Value = 12632.7495  # L122: 12965.727 -> 12632.7495
$ws.Cells.Item(122, 13).Value = -12411.799  # M122: -19937.125 -> -12411.799
$ws.Cells.Item(122, 14).Value = -17532.7495  # N122: -17865.727 -> -17532.7495

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(92, 8).Value = 107967  # H92: 96674.57000000001 -> 107967
$ws.Cells.Item(92, 10).Value = 107967  # J92: 96674.57000000001 -> 107967
$ws.Cells.Item(92, 12).Value = 107967  # L92: 96674.57000000001 -> 107967
$ws.Cells.Item(92, 14).Value = -112959  # N92: -101666.57 -> -112959

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 8584.286  # H4: 10000 -> 8584.286
$ws.Cells.Item(4, 9).Value = 100  # I4: 0 -> 100
$ws.Cells.Item(4, 10).Value = 9998.333000000001  # J4: 10000 -> 9998.333000000001
$ws.Cells.Item(4, 11).Value = 100  # K4: 0 -> 100
$ws.Cells.Item(4, 12).Value = 9998.333000000001  # L4: 10000 -> 9998.333000000001
$ws.Cells.Item(4, 13).Value = 12  # M4: None -> 12
$ws.Cells.Item(4, 14).Value = -10222.333  # N4: -10224 -> -10222.333

$ws.Cells.Item(5, 8).Value = 835.3077  # H5: 250.81818 -> 835.3077
$ws.Cells.Item(5, 9).Value = 241.66667  # I5: 270 -> 241.66667
$ws.Cells.Item(5, 10).Value = 1344.1428  # J5: 234.83333 -> 1344.1428
$ws.Cells.Item(5, 11).Value = 241.66667  # K5: 270 -> 241.66667
$ws.Cells.Item(5, 12).Value = 1344.1428  # L5: 234.83333 -> 1344.1428
$ws.Cells.Item(5, 13).Value = -129.66667  # M5: -158 -> -129.66667
$ws.Cells.Item(5, 14).Value = -1568.1428  # N5: -458.83333 -> -1568.1428

$ws.Cells.Item(8, 8).Value = 12528.143  # H8: 10963.25 -> 12528.143
$ws.Cells.Item(8, 9).Value = 0  # I8: 9 -> 0
$ws.Cells.Item(8, 11).Value = 0  # K8: 9 -> 0
$ws.Cells.Item(8, 13).ClearContents()  # M8: 131 -> (removed)

$ws.Cells.Item(75, 8).Value = 42000  # H75: 0 -> 42000
$ws.Cells.Item(75, 10).Value = 42000  # J75: 0 -> 42000
$ws.Cells.Item(75, 12).Value = 42000  # L75: 0 -> 42000
$ws.Cells.Item(75, 14).Value = -43996  # N75: None -> -43996

$ws.Cells.Item(78, 8).Value = 42000  # H78: 0 -> 42000
$ws.Cells.Item(78, 10).Value = 42000  # J78: 0 -> 42000
$ws.Cells.Item(78, 12).Value = 126000  # L78: 0 -> 126000
$ws.Cells.Item(78, 14).Value = -135984  # N78: None -> -135984

$ws.Cells.Item(86, 8).Value = 3140.0454  # H86: 8740.75 -> 3140.0454
$ws.Cells.Item(86, 9).Value = 2612.75  # I86: 13867.444 -> 2612.75
$ws.Cells.Item(86, 10).Value = 3772.8  # J86: 4546.1816 -> 3772.8
$ws.Cells.Item(86, 11).Value = 2612.75  # K86: 13867.444 -> 2612.75
$ws.Cells.Item(86, 12).Value = 3772.8  # L86: 4546.1816 -> 3772.8
$ws.Cells.Item(86, 13).Value = -1489.75  # M86: -12744.444 -> -1489.75
$ws.Cells.Item(86, 14).Value = -6018.8  # N86: -6792.1816 -> -6018.8

$ws.Cells.Item(87, 8).Value = 43997.5  # H87: 42865 -> 43997.5
$ws.Cells.Item(87, 9).Value = 44000  # I87: 0 -> 44000
$ws.Cells.Item(87, 10).Value = 43995  # J87: 42865 -> 43995
$ws.Cells.Item(87, 11).Value = 44000  # K87: 0 -> 44000
$ws.Cells.Item(87, 12).Value = 43995  # L87: 42865 -> 43995
$ws.Cells.Item(87, 13).Value = -42814  # M87: None -> -42814
$ws.Cells.Item(87, 14).Value = -46367  # N87: -45237 -> -46367

$ws.Cells.Item(89, 8).Value = 3140.0454  # H89: 8740.75 -> 3140.0454
$ws.Cells.Item(89, 9).Value = 2612.75  # I89: 13867.444 -> 2612.75
$ws.Cells.Item(89, 10).Value = 3772.8  # J89: 4546.1816 -> 3772.8
$ws.Cells.Item(89, 11).Value = 13063.75  # K89: 69337.22 -> 13063.75
$ws.Cells.Item(89, 12).Value = 18864  # L89: 22730.908 -> 18864
$ws.Cells.Item(89, 13).Value = -7447.75  # M89: -63721.22 -> -7447.75
$ws.Cells.Item(89, 14).Value = -30096  # N89: -33962.908 -> -30096

$ws.Cells.Item(90, 8).Value = 43997.5  # H90: 42865 -> 43997.5
$ws.Cells.Item(90, 9).Value = 44000  # I90: 0 -> 44000
$ws.Cells.Item(90, 10).Value = 43995  # J90: 42865 -> 43995
$ws.Cells.Item(90, 11).Value = 132000  # K90: 0 -> 132000
$ws.Cells.Item(90, 12).Value = 131985  # L90: 128595 -> 131985
$ws.Cells.Item(90, 13).Value = -126072  # M90: None -> -126072
$ws.Cells.Item(90, 14).Value = -143841  # N90: -140451 -> -143841

$ws.Cells.Item(94, 8).Value = 6593.0713  # H94: 6825.6294 -> 6593.0713
$ws.Cells.Item(94, 10).Value = 7651.609  # J94: 7985.136 -> 7651.609
$ws.Cells.Item(94, 12).Value = 7651.609  # L94: 7985.136 -> 7651.609
$ws.Cells.Item(94, 14).Value = -8553.609  # N94: -8887.136 -> -8553.609

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1966.4375  # H5: 1908.7059 -> 1966.4375
$ws.Cells.Item(5, 9).Value = 2332.5454  # I5: 2750.889 -> 2332.5454
$ws.Cells.Item(5, 10).Value = 1161  # J5: 961.25 -> 1161
$ws.Cells.Item(5, 11).Value = 6997.6362  # K5: 8252.667000000001 -> 6997.6362
$ws.Cells.Item(5, 12).Value = 3483  # L5: 2883.75 -> 3483
$ws.Cells.Item(5, 13).Value = -6885.6362  # M5: -8140.667000000001 -> -6885.6362
$ws.Cells.Item(5, 14).Value = -3707  # N5: -3107.75 -> -3707

$ws.Cells.Item(14, 8).Value = 58.615383  # H14: 30.916666 -> 58.615383
$ws.Cells.Item(14, 9).Value = 58.615383  # I14: 30.916666 -> 58.615383
$ws.Cells.Item(14, 11).Value = 175.846149  # K14: 92.74999800000001 -> 175.846149
$ws.Cells.Item(14, 13).Value = -2.846148999999997  # M14: 80.25000199999999 -> -2.846148999999997

$ws.Cells.Item(23, 8).Value = 274.0476  # H23: 248.40909 -> 274.0476
$ws.Cells.Item(23, 9).Value = 55.75  # I23: 46.6 -> 55.75
$ws.Cells.Item(23, 10).Value = 325.41177  # J23: 307.7647 -> 325.41177
$ws.Cells.Item(23, 11).Value = 167.25  # K23: 139.8 -> 167.25
$ws.Cells.Item(23, 12).Value = 976.23531  # L23: 923.2941000000001 -> 976.23531
$ws.Cells.Item(23, 13).Value = 67.75  # M23: 95.19999999999999 -> 67.75
$ws.Cells.Item(23, 14).Value = -1446.23531  # N23: -1393.2941 -> -1446.23531

$ws.Cells.Item(33, 8).Value = 1033.1111  # H33: 1124.8235 -> 1033.1111
$ws.Cells.Item(33, 9).Value = 807  # I33: 942 -> 807
$ws.Cells.Item(33, 10).Value = 1315.75  # J33: 1330.5 -> 1315.75
$ws.Cells.Item(33, 11).Value = 4842  # K33: 5652 -> 4842
$ws.Cells.Item(33, 12).Value = 7894.5  # L33: 7983 -> 7894.5
$ws.Cells.Item(33, 13).Value = -4559  # M33: -5369 -> -4559
$ws.Cells.Item(33, 14).Value = -8460.5  # N33: -8549 -> -8460.5

$ws.Cells.Item(107, 8).Value = 492.73685  # H107: 413.7857 -> 492.73685
$ws.Cells.Item(107, 9).Value = 522.3570999999999  # I107: 428.2381 -> 522.3570999999999
$ws.Cells.Item(107, 10).Value = 409.8  # J107: 370.42856 -> 409.8
$ws.Cells.Item(107, 11).Value = 1567.0713  # K107: 1284.7143 -> 1567.0713
$ws.Cells.Item(107, 12).Value = 1229.4  # L107: 1111.28568 -> 1229.4
$ws.Cells.Item(107, 13).Value = 352.9287000000002  # M107: 635.2857000000001 -> 352.9287000000002
$ws.Cells.Item(107, 14).Value = -5069.4  # N107: -4951.28568 -> -5069.4

$ws.Cells.Item(118, 8).Value = 2926.5454  # H118: 2828 -> 2926.5454
$ws.Cells.Item(118, 10).Value = 3079.2  # J118: 3032 -> 3079.2
$ws.Cells.Item(118, 12).Value = 9237.599999999999  # L118: 9096 -> 9237.599999999999
$ws.Cells.Item(118, 14).Value = -11723.6  # N118: -11582 -> -11723.6

$ws.Cells.Item(126, 8).Value = 4661.5386  # H126: 3795.625 -> 4661.5386
$ws.Cells.Item(126, 9).Value = 0  # I126: 3030 -> 0
$ws.Cells.Item(126, 10).Value = 4661.5386  # J126: 3846.6667 -> 4661.5386
$ws.Cells.Item(126, 11).Value = 0  # K126: 9090 -> 0
$ws.Cells.Item(126, 12).Value = 13984.6158  # L126: 11540.0001 -> 13984.6158
$ws.Cells.Item(126, 13).ClearContents()  # M126: -4150 -> (removed)
$ws.Cells.Item(126, 14).Value = -23864.6158  # N126: -21420.0001 -> -23864.6158

$ws.Cells.Item(129, 8).Value = 2632497.8  # H129: 3126094 -> 2632497.8
$ws.Cells.Item(129, 9).Value = 712.5  # I129: 780 -> 712.5
$ws.Cells.Item(129, 10).Value = 3334307.2  # J129: 3572567.2 -> 3334307.2
$ws.Cells.Item(129, 11).Value = 2137.5  # K129: 2340 -> 2137.5
$ws.Cells.Item(129, 12).Value = 10002921.6  # L129: 10717701.6 -> 10002921.6
$ws.Cells.Item(129, 13).Value = 2862.5  # M129: 2660 -> 2862.5
$ws.Cells.Item(129, 14).Value = -10012921.6  # N129: -10727701.6 -> -10012921.6

$ws.Cells.Item(131, 8).Value = 1397  # H131: 1420.2208 -> 1397
$ws.Cells.Item(131, 10).Value = 1502.2877  # J131: 1534.3043 -> 1502.2877
$ws.Cells.Item(131, 12).Value = 4506.8631  # L131: 4602.9129 -> 4506.8631
$ws.Cells.Item(131, 14).Value = -14586.8631  # N131: -14682.9129 -> -14586.8631

$ws.Cells.Item(132, 8).Value = 1985.8214  # H132: 1891.4333 -> 1985.8214
$ws.Cells.Item(132, 9).Value = 1580  # I132: 1488.8889 -> 1580
$ws.Cells.Item(132, 10).Value = 2613  # J132: 2495.25 -> 2613
$ws.Cells.Item(132, 11).Value = 14220  # K132: 13400.0001 -> 14220
$ws.Cells.Item(132, 12).Value = 23517  # L132: 22457.25 -> 23517
$ws.Cells.Item(132, 13).Value = -11690  # M132: -10870.0001 -> -11690
$ws.Cells.Item(132, 14).Value = -28577  # N132: -27517.25 -> -28577

$ws.Cells.Item(135, 8).Value = 1966.4375  # H135: 1908.7059 -> 1966.4375
$ws.Cells.Item(135, 9).Value = 2332.5454  # I135: 2750.889 -> 2332.5454
$ws.Cells.Item(135, 10).Value = 1161  # J135: 961.25 -> 1161
$ws.Cells.Item(135, 11).Value = 20992.9086  # K135: 24758.001 -> 20992.9086
$ws.Cells.Item(135, 12).Value = 10449  # L135: 8651.25 -> 10449
$ws.Cells.Item(135, 13).Value = -18457.9086  # M135: -22223.001 -> -18457.9086
$ws.Cells.Item(135, 14).Value = -15519  # N135: -13721.25 -> -15519

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(131, 8).Value = 33728  # H131: 33753 -> 33728
$ws.Cells.Item(131, 10).Value = 33728  # J131: 33753 -> 33728
$ws.Cells.Item(131, 12).Value = 33728  # L131: 33753 -> 33728
$ws.Cells.Item(131, 14).Value = -43808  # N131: -43833 -> -43808

$ws.Cells.Item(132, 8).Value = 1939.8572  # H132: 2297 -> 1939.8572
$ws.Cells.Item(132, 9).Value = 1645.1052  # I132: 1975.0209 -> 1645.1052
$ws.Cells.Item(132, 10).Value = 3232.2307  # J132: 3584.9167 -> 3232.2307
$ws.Cells.Item(132, 11).Value = 4935.3156  # K132: 5925.0627 -> 4935.3156
$ws.Cells.Item(132, 12).Value = 9696.6921  # L132: 10754.7501 -> 9696.6921
$ws.Cells.Item(132, 13).Value = -2405.3156  # M132: -3395.0627 -> -2405.3156
$ws.Cells.Item(132, 14).Value = -14756.6921  # N132: -15814.7501 -> -14756.6921

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(64, 8).Value = 34383.332  # H64: 35000 -> 34383.332
$ws.Cells.Item(64, 10).Value = 34383.332  # J64: 35000 -> 34383.332
$ws.Cells.Item(64, 12).Value = 34383.332  # L64: 35000 -> 34383.332
$ws.Cells.Item(64, 14).Value = -34833.332  # N64: -35450 -> -34833.332

$ws.Cells.Item(67, 8).Value = 34383.332  # H67: 35000 -> 34383.332
$ws.Cells.Item(67, 10).Value = 34383.332  # J67: 35000 -> 34383.332
$ws.Cells.Item(67, 12).Value = 34383.332  # L67: 35000 -> 34383.332
$ws.Cells.Item(67, 14).Value = -35943.332  # N67: -36560 -> -35943.332

$ws.Cells.Item(76, 8).Value = 18000  # H76: 50000 -> 18000
$ws.Cells.Item(76, 10).Value = 18000  # J76: 50000 -> 18000
$ws.Cells.Item(76, 12).Value = 18000  # L76: 50000 -> 18000
$ws.Cells.Item(76, 14).Value = -18676  # N76: -50676 -> -18676

$ws.Cells.Item(79, 8).Value = 18000  # H79: 50000 -> 18000
$ws.Cells.Item(79, 10).Value = 18000  # J79: 50000 -> 18000
$ws.Cells.Item(79, 12).Value = 18000  # L79: 50000 -> 18000
$ws.Cells.Item(79, 14).Value = -20340  # N79: -52340 -> -20340

$ws.Cells.Item(99, 8).Value = 228285.72  # H99: 201000 -> 228285.72
$ws.Cells.Item(99, 9).Value = 307800  # I99: 258166.67 -> 307800
$ws.Cells.Item(99, 11).Value = 307800  # K99: 258166.67 -> 307800
$ws.Cells.Item(99, 13).Value = -304805  # M99: -255171.67 -> -304805

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(108, 8).Value = 0  # H108: 23000 -> 0
$ws.Cells.Item(108, 10).Value = 0  # J108: 23000 -> 0
$ws.Cells.Item(108, 12).Value = 0  # L108: 23000 -> 0
$ws.Cells.Item(108, 14).ClearContents()  # N108: -30680 -> (removed)

$ws.Cells.Item(122, 8).Value = 12021533  # H122: 14207141 -> 12021533
$ws.Cells.Item(122, 9).Value = 16668476  # I122: 19232744 -> 16668476
$ws.Cells.Item(122, 10).Value = 5684791  # J122: 6947935 -> 5684791
$ws.Cells.Item(122, 11).Value = 50005428  # K122: 57698232 -> 50005428
$ws.Cells.Item(122, 12).Value = 17054373  # L122: 20843805 -> 17054373
$ws.Cells.Item(122, 13).Value = -50002978  # M122: -57695782 -> -50002978
$ws.Cells.Item(122, 14).Value = -17059273  # N122: -20848705 -> -17059273

$ws.Cells.Item(126, 8).Value = 6633  # H126: 10999 -> 6633
$ws.Cells.Item(126, 9).Value = 8953.200000000001  # I126: 10999 -> 8953.200000000001
$ws.Cells.Item(126, 10).Value = 832.5  # J126: 0 -> 832.5
$ws.Cells.Item(126, 11).Value = 26859.6  # K126: 32997 -> 26859.6
$ws.Cells.Item(126, 12).Value = 2497.5  # L126: 0 -> 2497.5
$ws.Cells.Item(126, 13).Value = -24389.6  # M126: -30527 -> -24389.6
$ws.Cells.Item(126, 14).Value = -7437.5  # N126: None -> -7437.5
